$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.060988187789917
$ws.Range("B1").Value = 2.640792846679688
$ws.Range("C1").Value = 2.825677633285522
$ws.Range("D1").Value = 3.398555755615234
$ws.Range("E1").Value = 2.256650447845459
